$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 26 (RM 232) entirely - causes rows below to shift up by one
$ws.Rows.Item(26).Delete()

# Delete what is now row 27 (was SC 92) entirely - causes rows below to shift up again
$ws.Rows.Item(27).Delete()

# At this point the former rows 29-35 have shifted up to rows 27-33,
# and former row 27 (SC 5) is now row 26. Re-randomize the "missing" pattern
# for column B across rows 27-32 (SC 101, SC 105, SC 119, SC 120, SC 193; SC132/SC232 unaffected).
$ws.Range("B27").Value = -20.4
$ws.Range("B28").ClearContents()
$ws.Range("B29").ClearContents()
$ws.Range("B30").Value = -19.7
$ws.Range("B32").ClearContents()
